$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated relevance marker time-to-discovery (TD) values: column C = td_sim_1,
# column D = record_atd (average of B and C), per row. Row 106 holds the
# average_simulation_TD summary and only its C value (average of C2:C105) changes.
$data = @(
    @{Row=2; C=164; D=157},
    @{Row=3; C=56; D=47},
    @{Row=4; C=63; D=78.5},
    @{Row=5; C=203; D=203},
    @{Row=6; C=44; D=49},
    @{Row=7; C=71; D=60.5},
    @{Row=8; C=132; D=124},
    @{Row=9; C=98; D=100.5},
    @{Row=10; C=134; D=137.5},
    @{Row=11; C=197; D=193},
    @{Row=12; C=162; D=166},
    @{Row=13; C=382; D=380},
    @{Row=14; C=123; D=128.5},
    @{Row=15; C=174; D=185.5},
    @{Row=16; C=89; D=90.5},
    @{Row=17; C=30; D=58},
    @{Row=18; C=93; D=83},
    @{Row=19; C=106; D=108.5},
    @{Row=20; C=136; D=108.5},
    @{Row=21; C=135; D=139},
    @{Row=22; C=428; D=435},
    @{Row=23; C=208; D=213.5},
    @{Row=24; C=176; D=176},
    @{Row=25; C=24; D=36},
    @{Row=26; C=13; D=11},
    @{Row=27; C=14; D=10},
    @{Row=28; C=31; D=74.5},
    @{Row=29; C=177; D=185},
    @{Row=30; C=54; D=45.5},
    @{Row=31; C=25; D=20.5},
    @{Row=32; C=273; D=273},
    @{Row=33; C=84; D=76.5},
    @{Row=34; C=51; D=42.5},
    @{Row=35; C=58; D=51.5},
    @{Row=36; C=74; D=67},
    @{Row=37; C=20; D=57.5},
    @{Row=38; C=42; D=37},
    @{Row=39; C=52; D=46.5},
    @{Row=40; C=97; D=97.5},
    @{Row=41; C=18; D=21},
    @{Row=42; C=103; D=115},
    @{Row=43; C=180; D=178.5},
    @{Row=44; C=65; D=59},
    @{Row=45; C=36; D=31},
    @{Row=46; C=81; D=76.5},
    @{Row=47; C=9; D=12},
    @{Row=48; C=95; D=101},
    @{Row=49; C=17; D=20},
    @{Row=50; C=50; D=46},
    @{Row=51; C=53; D=44},
    @{Row=52; C=113; D=108.5},
    @{Row=53; C=109; D=123.5},
    @{Row=54; C=161; D=170.5},
    @{Row=55; C=388; D=396.5},
    @{Row=56; C=171; D=163.5},
    @{Row=57; C=49; D=42.5},
    @{Row=58; C=7; D=5.5},
    @{Row=59; C=76; D=66.5},
    @{Row=60; C=62; D=52.5},
    @{Row=61; C=101; D=89.5},
    @{Row=62; C=46; D=39.5},
    @{Row=63; C=61; D=54},
    @{Row=64; C=114; D=117},
    @{Row=65; C=3; D=2},
    @{Row=66; C=128; D=128},
    @{Row=67; C=190; D=180.5},
    @{Row=68; C=183; D=175.5},
    @{Row=69; C=47; D=56.5},
    @{Row=70; C=37; D=32},
    @{Row=71; C=4; D=4},
    @{Row=72; C=34; D=29.5},
    @{Row=73; C=1767; D=1708},
    @{Row=74; C=91; D=92},
    @{Row=75; C=116; D=109},
    @{Row=76; C=129; D=129},
    @{Row=77; C=175; D=170},
    @{Row=78; C=38; D=33},
    @{Row=79; C=112; D=118},
    @{Row=80; C=12; D=10},
    @{Row=81; C=232; D=215.5},
    @{Row=82; C=187; D=185.5},
    @{Row=83; C=117; D=114.5},
    @{Row=84; C=55; D=47},
    @{Row=85; C=125; D=129},
    @{Row=86; C=140; D=138.5},
    @{Row=87; C=172; D=172.5},
    @{Row=88; C=149; D=144.5},
    @{Row=89; C=289; D=294.5},
    @{Row=90; C=16; D=14},
    @{Row=91; C=121; D=117},
    @{Row=92; C=19; D=58},
    @{Row=93; C=147; D=151},
    @{Row=94; C=60; D=54.5},
    @{Row=95; C=182; D=181.5},
    @{Row=96; C=122; D=122},
    @{Row=97; C=191; D=188},
    @{Row=98; C=94; D=91},
    @{Row=99; C=5; D=3.5},
    @{Row=100; C=39; D=34},
    @{Row=101; C=35; D=51},
    @{Row=102; C=216; D=212},
    @{Row=103; C=181; D=175},
    @{Row=104; C=165; D=199},
    @{Row=105; C=11; D=10.5},
    @{Row=106; C=123.0192307692308; D=$null}

)

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
    if ($null -ne $item.D) {
        $ws.Cells.Item($item.Row, 4).Value = $item.D
    }
}
